# Auto-generated PowerShell Excel COM-interop script
# Applies the 'cryptos list' data refresh described in the commit/diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. '26.769.76', '1.0000')
# that must stay literal text, matching the source data feed's formatting.
# Temporarily force Text number format so Excel doesn't auto-convert these
# strings into numbers, then clear the format back to the default.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.769.76'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.799.18'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '308.66'
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4596'
$ws.Range("E7").Value = '  +2.63%  '
$ws.Range("D8").Value = '0.3718'
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").Value = '0.07253'
$ws.Range("E9").Value = '  -3.49%  '
$ws.Range("D10").Value = '0.8560'
$ws.Range("E10").Value = '  -4.41%  '
$ws.Range("D11").Value = '20.39'
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("D12").Value = '1.791.55'
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '5.314'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '6.495'
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '0.07029'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").Value = '90.56'
$ws.Range("E16").Value = '  -4.40%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '0.000008617'
$ws.Range("E18").Value = '  -2.46%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '14.62'
$ws.Range("E20").Value = '  -4.01%  '
$ws.Range("D21").Value = '26.778.53'
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '10.64'
$ws.Range("E23").Value = '  -2.94%  '
$ws.Range("D24").Value = '2.032.07'
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("E25").Value = '  -4.83%  '
$ws.Range("D26").Value = '149.61'
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.149'
$ws.Range("E27").Value = '  -14.04%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.17'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = '5.217'
$ws.Range("E29").Value = '  -3.02%  '
$ws.Range("D30").Value = '114.01'
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("D31").Value = '0.08853'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '0.7542'
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("D33").Value = '1.159'
$ws.Range("E33").Value = '  -3.79%  '
$ws.Range("D34").Value = '4.428'
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("D35").Value = '2.885'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").Value = '0.01938'
$ws.Range("E38").Value = '  -2.77%  '
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.368'
$ws.Range("E40").Value = '  +3.61%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.896'
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '7.159'
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("D43").Value = '0.5229'
$ws.Range("E43").Value = '  -1.97%  '
$ws.Range("D44").Value = '0.1645'
$ws.Range("E44").Value = '  -5.15%  '
$ws.Range("D45").Value = '8.497'
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").Value = '0.4998'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").Value = '10.27'
$ws.Range("E47").Value = '  -4.53%  '
$ws.Range("D48").Value = '104.06'
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("D49").Value = '0.9999'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '1.643'
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("E51").Value = '  -1.19%  '

# Restore default (General) formatting on the Price cells so only the
# values changed, not the cell styling.
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
